$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: new TODO item - let user exit the server.
$ws.Range("K8").Value = "Make user able to exit server"
$ws.Range("L8").Value = "By hitting 'esc'"
$ws.Range("M8").Value = "TODO"
$ws.Range("M8").Interior.Color = $ws.Range("M4").Interior.Color

# Row 5: Networks section - "Get server to send maze back" is now underway,
# with a note on how it connects to receiving messages on the client.
$ws.Range("L5").Value = "Get client to receive messages"
$ws.Range("M5").Value = "UNDERWAY"
$ws.Range("M5").Interior.Color = $ws.Range("H12").Interior.Color

$ws.Range("M19").Select()
